$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2688.8
$ws.Range("I31").Value = 2688.8
$ws.Range("K31").Value = 8066.400000000001
$ws.Range("M31").Value = -7836.400000000001
$ws.Range("H40").Value = 1800
$ws.Range("J40").Value = 1800
$ws.Range("L40").Value = 1800
$ws.Range("N40").Value = -2150
$ws.Range("H44").Value = 17500
$ws.Range("J44").Value = 17500
$ws.Range("L44").Value = 17500
$ws.Range("N44").Value = -18424
$ws.Range("H74").Value = 17999.666
$ws.Range("I74").Value = 13999
$ws.Range("K74").Value = 13999
$ws.Range("M74").Value = -13063
$ws.Range("H77").Value = 17999.666
$ws.Range("I77").Value = 13999
$ws.Range("K77").Value = 69995
$ws.Range("M77").Value = -65315
$ws.Range("H112").Value = 78664.21000000001
$ws.Range("I112").Value = 3644.5
$ws.Range("J112").Value = 91167.5
$ws.Range("K112").Value = 10933.5
$ws.Range("L112").Value = 273502.5
$ws.Range("M112").Value = -9825.5
$ws.Range("N112").Value = -275718.5
$ws.Range("H127").Value = 250002830
$ws.Range("J127").Value = 5108.5
$ws.Range("L127").Value = 15325.5
$ws.Range("N127").Value = -25245.5
$ws.Range("H135").Value = 963.73334
$ws.Range("I135").Value = 963.73334
$ws.Range("K135").Value = 8673.600060000001
$ws.Range("M135").Value = -6138.600060000001
$ws.Range("H137").Value = 2330.389
$ws.Range("I137").Value = 2218.318
$ws.Range("J137").Value = 2506.5
$ws.Range("K137").Value = 6654.954000000001
$ws.Range("L137").Value = 7519.5
$ws.Range("M137").Value = -4104.954000000001
$ws.Range("N137").Value = -12619.5
$ws.Range("H138").Value = 4545.94
$ws.Range("J138").Value = 4808.4155
$ws.Range("L138").Value = 14425.2465
$ws.Range("N138").Value = -24705.2465

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10171.75
$ws.Range("I45").Value = 15795.286
$ws.Range("K45").Value = 15795.286
$ws.Range("M45").Value = -15418.286
$ws.Range("H61").Value = 450317.62
$ws.Range("I61").Value = 4716.4287
$ws.Range("K61").Value = 4716.4287
$ws.Range("M61").Value = -4504.4287
$ws.Range("H97").Value = 1642.5714
$ws.Range("I97").Value = 2059.8572
$ws.Range("K97").Value = 2059.8572
$ws.Range("M97").Value = -1563.8572
$ws.Range("H132").Value = 3245.8
$ws.Range("I132").Value = 3083.842
$ws.Range("J132").Value = 3525.5454
$ws.Range("K132").Value = 9251.526
$ws.Range("L132").Value = 10576.6362
$ws.Range("M132").Value = -6721.526
$ws.Range("N132").Value = -15636.6362
$ws.Range("H136").Value = 450317.62
$ws.Range("I136").Value = 4716.4287
$ws.Range("K136").Value = 14149.2861
$ws.Range("M136").Value = -11599.2861

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 11969.167
$ws.Range("I26").Value = 11969.167
$ws.Range("K26").Value = 11969.167
$ws.Range("M26").Value = -11677.167
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H96").Value = 15157.223
$ws.Range("I96").Value = 9651.714
$ws.Range("K96").Value = 9651.714
$ws.Range("M96").Value = -6905.714
$ws.Range("H99").Value = 5395.6
$ws.Range("I99").Value = 3524.2942
$ws.Range("K99").Value = 3524.2942
$ws.Range("M99").Value = -2026.2942
$ws.Range("H134").Value = 1298.4773
$ws.Range("I134").Value = 1146.0238
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 3438.0714
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -903.0713999999998
$ws.Range("N134").Value = -18570

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41144.96
$ws.Range("I31").Value = 47767.184
$ws.Range("J31").Value = 4722.75
$ws.Range("K31").Value = 47767.184
$ws.Range("L31").Value = 4722.75
$ws.Range("M31").Value = -47472.184
$ws.Range("N31").Value = -5312.75
$ws.Range("H34").Value = 41144.96
$ws.Range("I34").Value = 47767.184
$ws.Range("J34").Value = 4722.75
$ws.Range("K34").Value = 47767.184
$ws.Range("L34").Value = 4722.75
$ws.Range("M34").Value = -47565.184
$ws.Range("N34").Value = -5126.75
$ws.Range("H138").Value = 116666.5
$ws.Range("J138").Value = 117000
$ws.Range("L138").Value = 117000
$ws.Range("N138").Value = -127280
$ws.Range("H141").Value = 248420.86
$ws.Range("J141").Value = 261665.42
$ws.Range("L141").Value = 261665.42
$ws.Range("N141").Value = -272025.42

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 35768.93
$ws.Range("I2").Value = 41724.668
$ws.Range("K2").Value = 250348.008
$ws.Range("M2").Value = -250235.008
$ws.Range("H12").Value = 303.2
$ws.Range("J12").Value = 390.17648
$ws.Range("L12").Value = 1170.52944
$ws.Range("N12").Value = -1516.52944
$ws.Range("H45").Value = 25002208
$ws.Range("J45").Value = 2944
$ws.Range("L45").Value = 8832
$ws.Range("N45").Value = -9896
$ws.Range("H97").Value = 245.28572
$ws.Range("I97").Value = 145
$ws.Range("K97").Value = 435
$ws.Range("M97").Value = 61
$ws.Range("H99").Value = 3199.5
$ws.Range("I99").Value = 2400
$ws.Range("K99").Value = 7200
$ws.Range("M99").Value = -4954

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 87501
$ws.Range("I12").Value = 99998
$ws.Range("K12").Value = 99998
$ws.Range("M12").Value = -99858
$ws.Range("H35").Value = 234500
$ws.Range("I35").Value = 444000
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 444000
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -443702
$ws.Range("N35").Value = -25596
$ws.Range("H132").Value = 3934.375
$ws.Range("I132").Value = 2760.3333
$ws.Range("K132").Value = 8280.999899999999
$ws.Range("M132").Value = -5750.999899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 16559.25
$ws.Range("J43").Value = 16559.25
$ws.Range("L43").Value = 16559.25
$ws.Range("N43").Value = -16945.25
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2022.6119
$ws.Range("I132").Value = 1739.2881
$ws.Range("K132").Value = 5217.8643
$ws.Range("M132").Value = -2687.8643
$ws.Range("H135").Value = 109147.836
$ws.Range("J135").Value = 109147.836
$ws.Range("L135").Value = 109147.836
$ws.Range("N135").Value = -119287.836
$ws.Range("H136").Value = 2678.3333
$ws.Range("I136").Value = 2297.8696
$ws.Range("K136").Value = 6893.6088
$ws.Range("M136").Value = -4343.6088
